$wb = $excel.ActiveWorkbook

# Troponin sheet: move the selection over to E15
$troponin = $wb.Worksheets.Item("Troponin")
[void]$troponin.Activate()
[void]$troponin.Range("E15").Select()

# --- Finish the docking run for Anna: add a new results sheet -------------
# The new sheet reuses the same layout/styling as the "GO" sheet, so copy
# that sheet to the end of the workbook and then overwrite its data with
# the freshly completed docking numbers.
$goSheet = $wb.Worksheets.Item("GO")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$goSheet.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Sheet1"

# Row 5 - Dihidroquercetin
$newSheet.Range("D5").Value = -6.47
$newSheet.Range("E5").Value = 18.07
$newSheet.Range("F5").Value = -8.26
$newSheet.Range("G5").Value = -8.2
$newSheet.Range("H5").Value = -0.06
$newSheet.Range("I5").Value = -2.66
$newSheet.Range("J5").Value = 1.79
$newSheet.Range("K5").Value = -2.66

# Row 6 - Firulic acid
$newSheet.Range("D6").Value = -5.23
$newSheet.Range("E6").Value = 145.49
$newSheet.Range("F6").Value = -6.73
$newSheet.Range("G6").Value = -6.53
$newSheet.Range("H6").Value = -0.2
$newSheet.Range("I6").Value = -0.94
$newSheet.Range("J6").Value = 1.49
$newSheet.Range("K6").Value = -0.94

# Row 7 - Galloic acid
$newSheet.Range("D7").Value = -4.96
$newSheet.Range("E7").Value = 232.56
$newSheet.Range("F7").Value = -6.45
$newSheet.Range("G7").Value = -6.05
$newSheet.Range("H7").Value = -0.4
$newSheet.Range("I7").Value = -1.29
$newSheet.Range("J7").Value = 1.49
$newSheet.Range("K7").Value = -1.29

# Row 8 - Quercetin
$newSheet.Range("D8").Value = -6.39
$newSheet.Range("E8").Value = 20.61
$newSheet.Range("F8").Value = -8.18
$newSheet.Range("G8").Value = -7.6
$newSheet.Range("H8").Value = -0.58
$newSheet.Range("I8").Value = -2.47
$newSheet.Range("J8").Value = 1.79
$newSheet.Range("K8").Value = -2.47

# Row 9 - Rosavin
$newSheet.Range("D9").Value = -7.13
$newSheet.Range("E9").Value = 5.93
$newSheet.Range("F9").Value = -11.01
$newSheet.Range("G9").Value = -10.68
$newSheet.Range("H9").Value = -0.33
$newSheet.Range("I9").Value = -3.88
$newSheet.Range("J9").Value = 3.88
$newSheet.Range("K9").Value = -3.88

# Row 10 - Rutin (Ki reported as a concentration string, not a number)
$newSheet.Range("D10").Value = -3.78
$newSheet.Range("E10").Value = "1.69 mM"
$newSheet.Range("F10").Value = -8.55
$newSheet.Range("G10").Value = -8.54
$newSheet.Range("H10").Value = -0.02
$newSheet.Range("I10").Value = -8.48
$newSheet.Range("J10").Value = 4.77
$newSheet.Range("K10").Value = -8.48

# Row 11 - Salidrodside
$newSheet.Range("D11").Value = -5.26
$newSheet.Range("E11").Value = 138.99
$newSheet.Range("F11").Value = -8.24
$newSheet.Range("G11").Value = -8.08
$newSheet.Range("H11").Value = -0.16
$newSheet.Range("I11").Value = -1.63
$newSheet.Range("J11").Value = 2.98
$newSheet.Range("K11").Value = -1.63

[void]$newSheet.Range("D5").Select()
